# Updates crypto price/volume figures per the authoritative source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.398.78'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.800.87'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '''225.16'
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").Value = '''0.599'
$ws.Range("E6").Value = '  +3.31%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = '''36.23'
$ws.Range("E8").Value = '  +3.39%  '
$ws.Range("E9").Value = '  -2.30%  '
$ws.Range("D10").Value = '''0.0678'
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("D11").Value = '''0.0965'
$ws.Range("E11").Value = '  +1.36%  '
$ws.Range("D12").Value = '2.061.20'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '''11.27'
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").Value = '1.832.50'
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D16").Value = '34.389.58'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '''4.43'
$ws.Range("E17").Value = '  +2.50%  '
$ws.Range("D18").Value = '''68.49'
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").Value = '''242.27'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = '0.0₃0772'
$ws.Range("E20").Value = '  -2.83%  '
$ws.Range("D21").Value = '''11.25'
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("E24").Value = '  +5.07%  '
$ws.Range("D25").Value = '''170.70'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("E26").Value = '  +4.40%  '
$ws.Range("D27").Value = '''17.35'
$ws.Range("E27").Value = '  +3.95%  '
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").Value = '''3.91'
$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("E32").Value = '  -1.13%  '
$ws.Range("D33").Value = '''0.0514'
$ws.Range("E33").Value = '  -2.37%  '
$ws.Range("E34").Value = '  -3.12%  '
$ws.Range("D35").Value = '1.362.09'
$ws.Range("E35").Value = '  -2.57%  '
$ws.Range("E36").Value = '  -4.06%  '
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("E38").Value = '  -6.38%  '
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("E41").Value = '  -1.82%  '
$ws.Range("D42").Value = '''80.68'
$ws.Range("E42").Value = '  -2.61%  '
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("E44").Value = '  +5.08%  '
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("D46").Value = '''0.0499'
$ws.Range("E46").Value = '  -2.69%  '
$ws.Range("D47").Value = '1.964.40'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("E48").Value = '  -3.53%  '
$ws.Range("E49").Value = '  -0.24%  '
$ws.Range("D50").Value = '''102.27'
$ws.Range("E50").Value = '  -1.95%  '
$ws.Range("E51").Value = '  -3.24%  '
